# Refactor goal status table to be more descriptive
#
# The template previously assumed the goal-status table was the only
# dynamically generated table in the document, so its Jinja loop variable
# was named generically ("tbl_contents"). This renames it to the more
# descriptive "goal_status_table" and, along the way, cleans up the
# document so the surrounding merge-field / template-tag text runs in
# each affected paragraph are consolidated into single runs (removing
# stale spell-check markup left over from how those runs were typed).

$d = $word.ActiveDocument

function Retype-Text([string]$oldText, [string]$newText) {
    # Re-typing a run's text (even to the same value) makes Word
    # re-flow/merge the run(s) spanning that text and drop any stale
    # proofing-error markers. We first swap in a disposable placeholder
    # (guaranteed not to equal $oldText) and then set the real value, so
    # this also works when $newText -eq $oldText.
    $rng = $d.Content
    $ok = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $oldText"
    }
    $rng.Text = "@@TEMP_PLACEHOLDER@@"

    $rng2 = $d.Content
    $ok2 = $rng2.Find.Execute("@@TEMP_PLACEHOLDER@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok2) {
        throw "Could not find placeholder while retyping: $oldText"
    }
    $rng2.Text = $newText
}

# 1) "Goal Status across {{" / "agency_name" / "}}" -> single run
Retype-Text "Goal Status across {{agency_name}}" "Goal Status across {{agency_name}}"

# 2) "{{" / "previous_quarter_and_year" / "}}" -> single run
Retype-Text "{{previous_quarter_and_year}}" "{{previous_quarter_and_year}}"

# 3) "{{" / "current_quarter_and_year" / "}}" -> single run
Retype-Text "{{current_quarter_and_year}}" "{{current_quarter_and_year}}"

# 4) Rename the goal-status table loop variable: tbl_contents -> goal_status_table
#    "{%tr for item in " / "tbl_contents" -> "{%tr for item " / "in goal_status_table"
$rngA = $d.Content
$okA = $rngA.Find.Execute("item in ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okA) { throw "Could not find 'item in '" }
$rngA.Text = "item "

$rngB = $d.Content
$okB = $rngB.Find.Execute("tbl_contents", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okB) { throw "Could not find 'tbl_contents'" }
$rngB.Text = "in goal_status_table"

# 5) " {%" / "tc" / " for col in " / "item.cols" / " %}" -> single run
Retype-Text " {%tc for col in item.cols %}" " {%tc for col in item.cols %}"

# 6) " {%" / "tc" / " " / "endfor" / " %}" -> single run
Retype-Text " {%tc endfor %}" " {%tc endfor %}"

# 7) "{%tr " / "endfor" / " %}" -> single run (goal status table closing tag)
Retype-Text "{%tr endfor %}" "{%tr endfor %}"

# 8) " in {{ " / "current_quarter_and_year" / "}}" -> single run
Retype-Text " in {{ current_quarter_and_year}}" " in {{ current_quarter_and_year}}"

# 9) "{%tr for item in " / "challenge_count_table" / " %}" -> single run
Retype-Text "{%tr for item in challenge_count_table %}" "{%tr for item in challenge_count_table %}"

# 10) "{{ " / "item.col." / "count" / " }}" -> "{{ item.col." + "count" + " }}"
Retype-Text "{{ item.col." "{{ item.col."

# 11) "{%tr " / "endfor" / " %}" -> single run (challenge count table closing tag)
Retype-Text "{%tr endfor %}" "{%tr endfor %}"
